$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the formula-text values in column A. A leading apostrophe keeps
# the text quote-prefixed (same as the original author-entered cells),
# which preserves the existing quotePrefix cell style instead of Excel
# re-styling the cell as plain text.
$ws.Range("A3").Value = "'5*3-4+100"
$ws.Range("A4").Value = "'90+200+1000-245/10"

# Widen column A to 25 characters and drop the "best fit" auto-sizing that
# was previously in effect. (24.1666... is the COM ColumnWidth input that
# round-trips through Excel's pixel-based storage to an exact stored width
# of 25 - entering 25 directly overshoots to ~25.83 because of the
# character->pixel->character rounding Excel applies.)
$ws.Columns("A").ColumnWidth = 24.16666666666667

# Move the active selection to B4
$ws.Range("B4").Select()
